$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.896.39"
$ws.Range("E2").Value = "  +2.55%  "

$ws.Range("D3").Value = "3.593.03"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("E4").Value = "  -0.11%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "202.13"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +8.43%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "570.68"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "

$ws.Range("D7").Value = "3.591.25"
$ws.Range("E7").Value = "  +1.58%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.615"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.16%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.09%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.680"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +2.67%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "60.55"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +12.23%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.148"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000282"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +10.65%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "10.31"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.34%  "

$ws.Range("D15").Value = "4.155.40"
$ws.Range("E15").Value = "  +1.34%  "

$ws.Range("D16").Value = "3.593.17"
$ws.Range("E16").Value = "  +1.93%  "

$ws.Range("E17").Value = "  +1.17%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "19.13"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +5.71%  "

$ws.Range("D19").Value = "67.683.04"
$ws.Range("E19").Value = "  +2.19%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.28"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.30%  "

$ws.Range("E21").Value = "  +0.93%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "404.32"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +4.74%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "12.73"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +15.62%  "

$ws.Range("E24").Value = "  +0.81%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "84.77"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.89"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "3.90"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +11.06%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "12.45"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.22"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +4.91%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.93%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "31.54"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +2.51%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "676.20"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +9.63%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "12.13"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.01%  "

$ws.Range("E34").Value = "  -0.03%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "63.29"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.73%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "41.41"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.15%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.410"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.43%  "

$ws.Range("E38").Value = "  -0.07%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "3.21"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +11.19%  "

$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  +4.85%  "

$ws.Range("D41").Value = "3.200.63"
$ws.Range("E41").Value = "  +4.15%  "

$ws.Range("E42").Value = "  +1.11%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.70"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +7.30%  "

$ws.Range("E45").Value = "  +14.11%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.77"
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0410"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "

$ws.Range("E48").Value = "  +1.86%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "3.07"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -1.70%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "8.64"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.47%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "138.44"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.87%  "
